$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period values (Periodo Mora) in E16:E22 to ascending order 2405-2411
# (previously listed in descending order 2411-2405)
$ws.Range("E16").Value = "2405"
$ws.Range("E17").Value = "2406"
$ws.Range("E18").Value = "2407"
$ws.Range("E19").Value = "2408"
$ws.Range("E20").Value = "2409"
$ws.Range("E21").Value = "2410"
$ws.Range("E22").Value = "2411"

# Swap the Valor Mora amounts between the first and last rows of the table
$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 24266
